$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row changes
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Capitalize connector words (de/del/el/la/los/las) in state/municipality names
$ws.Range("B26").Value = "Marqués De Comillas"
$ws.Range("B27").Value = "Mazapa De Madero"
$ws.Range("B33").Value = "San Cristóbal De Las Casas"
$ws.Range("A49").Value = "Ciudad De México"
$ws.Range("A65").Value = "Estado De México"
$ws.Range("B67").Value = "Ecatepec De Morelos"
$ws.Range("B69").Value = "Naucalpan De Juárez"
$ws.Range("B74").Value = "Tlalnepantla De Baz"
$ws.Range("B77").Value = "Apaseo El Grande"
$ws.Range("B85").Value = "San Diego De La Unión"
$ws.Range("B87").Value = "San Luis De La Paz"
$ws.Range("B89").Value = "Acapulco De Juárez"
$ws.Range("B94").Value = "Atlamajalcingo Del Monte"
$ws.Range("B95").Value = "Ayutla De Los Libres"
$ws.Range("B97").Value = "Chilpancingo De Los Bravo"
$ws.Range("B105").Value = "Tlapa De Comonfort"
$ws.Range("B109").Value = "Huasca De Ocampo"
$ws.Range("B110").Value = "Pachuca De Soto"
$ws.Range("B111").Value = "Tulancingo De Bravo"
$ws.Range("B113").Value = "Autlán De Navarro"
$ws.Range("B116").Value = "Cuautitlán De García Barragán"
$ws.Range("B119").Value = "La Manzanilla De La Paz"
$ws.Range("B120").Value = "Lagos De Moreno"
$ws.Range("B124").Value = "Tizapán El Alto"
$ws.Range("B125").Value = "Tlajomulco De Zúñiga"
$ws.Range("B129").Value = "Unión De Tula"
$ws.Range("B131").Value = "Zapotlán El Grande"
$ws.Range("B153").Value = "Puente De Ixtla"
$ws.Range("B154").Value = "Tetela Del Volcán"
$ws.Range("B163").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B164").Value = "Oaxaca De Juárez"
$ws.Range("B168").Value = "San Dionisio Del Mar"
$ws.Range("B173").Value = "San Pedro El Alto"
$ws.Range("B192").Value = "Izúcar De Matamoros"
$ws.Range("B210").Value = "Tuzamapan De Galeana"
$ws.Range("B227").Value = "Villa De Ramos"
$ws.Range("B228").Value = "Villa De Reyes"
$ws.Range("B245").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B246").Value = "Muñoz De Domingo Arenas"
$ws.Range("B263").Value = "Lerdo De Tejada"
$ws.Range("B264").Value = "Martínez De La Torre"
$ws.Range("B272").Value = "Sayula De Alemán"

# Remove footer/metadata rows 284-288 (trailing notes below the data table)
$ws.Range("A284:A288").ClearContents()

Write-Output "done"
